$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.14506580613478093,
    -0.0059999999697843975,
    -0.0039999999743773884,
    -0.0079999999523838738,
    -0.0029999999745502492,
    -0.0019999999736057816,
    -0.0090916346680076643,
    -0.0099999999329010159,
    -0.0019999999716251438,
    -0.0019999999713409267,
    -0.0029999999663994359,
    -0.0034999999640286106,
    -0.0034999999647995494,
    -0.0079999999431796809,
    -0.00099999997841493382,
    0.024971880434306204,
    -0.0019999999739868102,
    -0.0039999999639555028,
    -0.0039999999795856667,
    -0.0039999999793369767,
    0.007986400478676714,
    -0.0039999999792286189,
    -0.0049999999678265183,
    -0.019999999891151532,
    -0.019999999889711795,
    -0.0024999999704391485,
    -0.0024999999694772512,
    0.018309321434324843,
    -0.0069999999390963907,
    -0.059999999676918936,
    0.023966062317695602,
    -0.0099999999226660918,
    -0.0039999999520095031
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
